# Update Leve profit-calculation figures (currentAveragePrice / LevePrice / LeveProfit
# columns H,I,J,K,L,M,N) across all Disciple-of-the-Hand sheets, per the scheduled
# market-price refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9123.857
$ws.Range("I74").Value = 8144.5
$ws.Range("K74").Value = 8144.5
$ws.Range("M74").Value = -7208.5
$ws.Range("H77").Value = 9123.857
$ws.Range("I77").Value = 8144.5
$ws.Range("K77").Value = 40722.5
$ws.Range("M77").Value = -36042.5
$ws.Range("H86").Value = 2911
$ws.Range("I86").Value = 3565.6667
$ws.Range("K86").Value = 3565.6667
$ws.Range("M86").Value = -2442.6667
$ws.Range("H89").Value = 2911
$ws.Range("I89").Value = 3565.6667
$ws.Range("K89").Value = 17828.3335
$ws.Range("M89").Value = -12212.3335
$ws.Range("H100").Value = 55723.582
$ws.Range("I100").Value = 65200.1
$ws.Range("J100").Value = 8341
$ws.Range("K100").Value = 65200.1
$ws.Range("L100").Value = 8341
$ws.Range("M100").Value = -64659.1
$ws.Range("N100").Value = -9423
$ws.Range("H132").Value = 3881632.2
$ws.Range("I132").Value = 4487753.5
$ws.Range("K132").Value = 13463260.5
$ws.Range("M132").Value = -13460730.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20292.145
$ws.Range("I32").Value = 21055.625
$ws.Range("J32").Value = 13166.333
$ws.Range("K32").Value = 21055.625
$ws.Range("L32").Value = 13166.333
$ws.Range("M32").Value = -20768.625
$ws.Range("N32").Value = -13740.333
$ws.Range("H61").Value = 4089.7058
$ws.Range("I61").Value = 1643.8334
$ws.Range("K61").Value = 1643.8334
$ws.Range("M61").Value = -1431.8334
$ws.Range("H62").Value = 54999
$ws.Range("J62").Value = 54999
$ws.Range("L62").Value = 54999
$ws.Range("N62").Value = -56247
$ws.Range("H65").Value = 54999
$ws.Range("J65").Value = 54999
$ws.Range("L65").Value = 164997
$ws.Range("N65").Value = -171237
$ws.Range("H119").Value = 80000
$ws.Range("J119").Value = 80000
$ws.Range("L119").Value = 80000
$ws.Range("N119").Value = -89676
$ws.Range("H121").Value = 80000
$ws.Range("J121").Value = 80000
$ws.Range("L121").Value = 80000
$ws.Range("N121").Value = -83494
$ws.Range("H122").Value = 1704.6
$ws.Range("I122").Value = 1602.4231
$ws.Range("K122").Value = 4807.2693
$ws.Range("M122").Value = -2357.2693
$ws.Range("H136").Value = 4089.7058
$ws.Range("I136").Value = 1643.8334
$ws.Range("K136").Value = 4931.5002
$ws.Range("M136").Value = -2381.5002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1453.75
$ws.Range("I86").Value = 1159.138
$ws.Range("K86").Value = 1159.138
$ws.Range("M86").Value = -36.13799999999992
$ws.Range("H89").Value = 1453.75
$ws.Range("I89").Value = 1159.138
$ws.Range("K89").Value = 5795.69
$ws.Range("M89").Value = -179.6899999999996
$ws.Range("H94").Value = 583.9737
$ws.Range("J94").Value = 910
$ws.Range("L94").Value = 910
$ws.Range("N94").Value = -1812
$ws.Range("H134").Value = 3290.361
$ws.Range("I134").Value = 3059.182
$ws.Range("K134").Value = 9177.545999999998
$ws.Range("M134").Value = -6642.545999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 23117.8
$ws.Range("I58").Value = 1597.6154
$ws.Range("K58").Value = 1597.6154
$ws.Range("M58").Value = -1394.6154
$ws.Range("H86").Value = 20543.889
$ws.Range("I86").Value = 11898.3
$ws.Range("K86").Value = 11898.3
$ws.Range("M86").Value = -10775.3
$ws.Range("H88").Value = 41666.668
$ws.Range("J88").Value = 41666.668
$ws.Range("L88").Value = 41666.668
$ws.Range("N88").Value = -42478.668
$ws.Range("H89").Value = 20543.889
$ws.Range("I89").Value = 11898.3
$ws.Range("K89").Value = 59491.5
$ws.Range("M89").Value = -53875.5
$ws.Range("H91").Value = 41666.668
$ws.Range("J91").Value = 41666.668
$ws.Range("L91").Value = 41666.668
$ws.Range("N91").Value = -44474.668
$ws.Range("H105").Value = 1347.85
$ws.Range("I105").Value = 902.05884
$ws.Range("K105").Value = 902.05884
$ws.Range("M105").Value = 844.94116
$ws.Range("H132").Value = 25994.291
$ws.Range("I132").Value = 29426.857
$ws.Range("K132").Value = 88280.571
$ws.Range("M132").Value = -85750.571
$ws.Range("H134").Value = 2736.8096
$ws.Range("I134").Value = 2263.4666
$ws.Range("K134").Value = 6790.399800000001
$ws.Range("M134").Value = -4255.399800000001
$ws.Range("H136").Value = 23117.8
$ws.Range("I136").Value = 1597.6154
$ws.Range("K136").Value = 4792.8462
$ws.Range("M136").Value = -2242.8462
$ws.Range("H141").Value = 163999.92
$ws.Range("J141").Value = 163999.92
$ws.Range("L141").Value = 163999.92
$ws.Range("N141").Value = -174359.92

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3140.5715
$ws.Range("J75").Value = 3373.5
$ws.Range("L75").Value = 10120.5
$ws.Range("N75").Value = -12116.5
$ws.Range("H78").Value = 3140.5715
$ws.Range("J78").Value = 3373.5
$ws.Range("L78").Value = 30361.5
$ws.Range("N78").Value = -40345.5
$ws.Range("H137").Value = 2556.5715
$ws.Range("I137").Value = 2407.077
$ws.Range("J137").Value = 4500
$ws.Range("K137").Value = 7221.231000000001
$ws.Range("L137").Value = 13500
$ws.Range("M137").Value = -2121.231000000001
$ws.Range("N137").Value = -23700

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 19314.834
$ws.Range("J15").Value = 19314.834
$ws.Range("L15").Value = 19314.834
$ws.Range("N15").Value = -19890.834
$ws.Range("H81").Value = 19314.834
$ws.Range("J81").Value = 19314.834
$ws.Range("L81").Value = 19314.834
$ws.Range("N81").Value = -21310.834
$ws.Range("H84").Value = 19314.834
$ws.Range("J84").Value = 19314.834
$ws.Range("L84").Value = 57944.50199999999
$ws.Range("N84").Value = -67928.50199999999
$ws.Range("H97").Value = 1256.4
$ws.Range("I97").Value = 1041.6111
$ws.Range("J97").Value = 1808.7142
$ws.Range("K97").Value = 1041.6111
$ws.Range("L97").Value = 1808.7142
$ws.Range("M97").Value = -545.6111000000001
$ws.Range("N97").Value = -2800.7142
$ws.Range("H102").Value = 2162.24
$ws.Range("I102").Value = 1872.1305
$ws.Range("K102").Value = 1872.1305
$ws.Range("M102").Value = -250.1305
$ws.Range("H122").Value = 3038.2222
$ws.Range("J122").Value = 2786.7144
$ws.Range("L122").Value = 8360.143199999999
$ws.Range("N122").Value = -13260.1432
$ws.Range("H132").Value = 6962.5
$ws.Range("I132").Value = 6962.5
$ws.Range("K132").Value = 20887.5
$ws.Range("M132").Value = -18357.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3403.3
$ws.Range("I40").Value = 3076.6191
$ws.Range("J40").Value = 4165.5557
$ws.Range("K40").Value = 3076.6191
$ws.Range("L40").Value = 4165.5557
$ws.Range("M40").Value = -2940.6191
$ws.Range("N40").Value = -4437.5557
$ws.Range("H82").Value = 9066.666999999999
$ws.Range("J82").Value = 11300
$ws.Range("L82").Value = 11300
$ws.Range("N82").Value = -12022
$ws.Range("H85").Value = 9066.666999999999
$ws.Range("J85").Value = 11300
$ws.Range("L85").Value = 11300
$ws.Range("N85").Value = -13796
$ws.Range("H93").Value = 3568.9092
$ws.Range("I93").Value = 3212.8235
$ws.Range("K93").Value = 3212.8235
$ws.Range("M93").Value = -1964.8235
$ws.Range("H119").Value = 132999
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 132999
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 132999
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -142675
$ws.Range("H132").Value = 7879.4
$ws.Range("I132").Value = 8000
$ws.Range("K132").Value = 24000
$ws.Range("M132").Value = -21470

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 102000
$ws.Range("I51").Value = 74000
$ws.Range("J51").Value = 130000
$ws.Range("K51").Value = 74000
$ws.Range("L51").Value = 130000
$ws.Range("M51").Value = -73490
$ws.Range("N51").Value = -131020
$ws.Range("H54").Value = 25702.666
$ws.Range("J54").Value = 31135
$ws.Range("L54").Value = 31135
$ws.Range("N54").Value = -32175
$ws.Range("H75").Value = 78688.75
$ws.Range("J75").Value = 78688.75
$ws.Range("L75").Value = 78688.75
$ws.Range("N75").Value = -80560.75
$ws.Range("H78").Value = 78688.75
$ws.Range("J78").Value = 78688.75
$ws.Range("L78").Value = 236066.25
$ws.Range("N78").Value = -245426.25
$ws.Range("H100").Value = 1171.5834
$ws.Range("I100").Value = 806
$ws.Range("K100").Value = 1612
$ws.Range("M100").Value = -1071
$ws.Range("H132").Value = 20331.95
$ws.Range("I132").Value = 25375.451
$ws.Range("K132").Value = 76126.353
$ws.Range("M132").Value = -73596.353
